$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("100Mbps 25ms")

# ---------------------------------------------------------------------------
# 1. New benchmark column E: "kcptun -fast Performance (MBytes/sec)"
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "kcptun -fast Performance (MBytes/sec)"
$ws.Range("E2").Value = 9.1199999999999992
$ws.Range("E3").Value = 8.25
$ws.Range("E4").Value = 8
$ws.Range("E5").Value = 6.75
$ws.Range("E6").Value = 4.0999999999999996

# Match column E's width to the other data columns (closest value this
# engine's pixel-quantized ColumnWidth model can reach to the authored 33.5).
$ws.Columns.Item(5).ColumnWidth = 32.785714285714285

# ---------------------------------------------------------------------------
# 2. Add the new series to the bar chart and slot it in right before the
#    Dragonite series (so it plots between TCP-BBR and Dragonite).
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$newSeries = $chart.SeriesCollection().NewSeries()
$newSeries.Values = '=''100Mbps 25ms''!$E$2:$E$6'
$newSeries.XValues = '=''100Mbps 25ms''!$A$2:$A$6'
$newSeries.Name = '=''100Mbps 25ms''!$E$1'

# The new series was appended at the end (plot order = last); move it so it
# plots immediately before the Dragonite series, matching the authored order.
$lastIndex = $chart.SeriesCollection().Count
$chart.SeriesCollection().Item($lastIndex).PlotOrder = $lastIndex - 1

# ---------------------------------------------------------------------------
# 3. Resize the chart's anchor now that the plotted data grew (moved down /
#    right by roughly a row+column of slack) — match the authored extent.
# ---------------------------------------------------------------------------
$co.Width = 474.750039370079
$co.Height = 262.874881889764

# ---------------------------------------------------------------------------
# 4. Selection moves to D23 in the saved view state.
# ---------------------------------------------------------------------------
[void]$ws.Range("D23").Select()
